$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: issue number and report week-covering dates ---
$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# --- Column E width adjustment (auto best-fit widened for new "500" value) ---
$ws.Columns("E").ColumnWidth = 7.433768

# --- Weekly crime-statistics table updates (rows 15-28, 33) ---
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 500
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = 42.857142857142
$ws.Range("L15").Value = 11.111111111111
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -56.521739130434
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 75
$ws.Range("I16").Value = 53
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = 1.923076923076
$ws.Range("L16").Value = 39.473684210526
$ws.Range("M16").Value = -32.051282051282
$ws.Range("N16").Value = -87.558685446009
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 5
$ws.Range("H17").Value = -44.444444444444
$ws.Range("I17").Value = 91
$ws.Range("J17").Value = 108
$ws.Range("K17").Value = -15.74074074074
$ws.Range("L17").Value = 8.333333333333
$ws.Range("M17").Value = 16.666666666666
$ws.Range("N17").Value = -59.375
$ws.Range("C18").Value = "0"
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("J18").Value = 78
$ws.Range("K18").Value = -8.974358974358
$ws.Range("L18").Value = -6.578947368421
$ws.Range("M18").Value = -63.020833333333
$ws.Range("N18").Value = -91.537544696066
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -12.5
$ws.Range("I19").Value = 279
$ws.Range("J19").Value = 334
$ws.Range("K19").Value = -16.467065868263
$ws.Range("L19").Value = -32.445520581113
$ws.Range("M19").Value = 14.344262295082
$ws.Range("N19").Value = -20.51282051282
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -13.333333333333
$ws.Range("I20").Value = 155
$ws.Range("J20").Value = 98
$ws.Range("K20").Value = 58.163265306122
$ws.Range("L20").Value = 82.35294117647
$ws.Range("M20").Value = 27.049180327868
$ws.Range("N20").Value = -89.55525606469
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -5.263157894736
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = -1.428571428571
$ws.Range("I21").Value = 660
$ws.Range("J21").Value = 679
$ws.Range("K21").Value = -2.798232695139
$ws.Range("L21").Value = -6.647807637906
$ws.Range("M21").Value = -8.460471567267
$ws.Range("N21").Value = -80.304386750223
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -18.181818181818
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -14.516129032258
$ws.Range("I24").Value = 1114
$ws.Range("J24").Value = 1178
$ws.Range("K24").Value = -5.432937181663
$ws.Range("L24").Value = -23.172413793103
$ws.Range("M24").Value = 21.350762527233
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -19.047619047619
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = -31.25
$ws.Range("I25").Value = 690
$ws.Range("J25").Value = 617
$ws.Range("K25").Value = 11.831442463533
$ws.Range("L25").Value = -21.142857142857
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 24
$ws.Range("I26").Value = 278
$ws.Range("J26").Value = 286
$ws.Range("K26").Value = -2.797202797202
$ws.Range("L26").Value = 18.297872340425
$ws.Range("M26").Value = -6.397306397306
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 250
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = 27.272727272727
$ws.Range("L27").Value = 16.666666666666
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("L33").Value = 350
